$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "Thống kê môn"

$ws.Range("A1").Value = "Năm"
$ws.Range("B1").Value = "Khu vực"
$ws.Range("C1").Value = "Tỉnh/Thành phố"
$ws.Range("D1").Value = "Môn"
$ws.Range("E1").Value = "Số lượng thí sinh"
$ws.Range("F1").Value = "Điểm trung bình"
$ws.Range("G1").Value = "Số lượng điểm < 1"
$ws.Range("H1").Value = "Số lượng điểm < 5"
$ws.Range("I1").Value = "Số lượng điểm >= 9"
$ws.Range("J1").Value = "Điểm cao nhất"
$ws.Range("K1").Value = "Điểm thấp nhất"
$ws.Range("L1").Value = "Kinh độ"
$ws.Range("M1").Value = "Vĩ độ"
